$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The MN now controls (almost) every system metadata field; only the
# "replica" row keeps CN as its controller. Update the ControlledBy (C),
# Method (E) and ModifiableBy (F) columns accordingly.

# serialVersion: ControlledBy CN -> MN, ModifiableBy CN-service-subject -> MN-service-subject
$ws.Range("C2").Value = "MN"
$ws.Range("F2").Value = "MN-service-subject"

# rightsHolder: Method CNAuthorization.setOwner -> MNAuthorization.setRightsHolder()
$ws.Range("E8").Value = "MNAuthorization.setRightsHolder()"

# accessPolicy: Method updated to the manual/Tier 1 MN process
$ws.Range("E9").Value = "manual (Tier 1), MNAuthorization.setAccessPolicy(), MNStorage.update ()(all must call CNAuthorization.systemMetadataChanged())"
$ws.Rows.Item(9).RowHeight = 45

# obsoletes: Method MNStorage.update -> MNStorage.update()
$ws.Range("E11").Value = "MNStorage.update()"

# obsoletedBy: Method MNStorage.update -> MNCore.setObsoletedBy(), MNStorage.update()
$ws.Range("E12").Value = "MNCore.setObsoletedBy(), MNStorage.update()"

# archived: Method MNStorage.delete -> MNCore.archive()
$ws.Range("E13").Value = "MNCore.archive()"

# dateSysMetadataModified: ControlledBy CN -> MN
$ws.Range("C15").Value = "MN"

# originMemberNode: ControlledBy CN -> MN
$ws.Range("C16").Value = "MN"

# authoritativeMemberNode: ControlledBy CN -> MN/CN, ModifiableBy updated
$ws.Range("C17").Value = "MN/CN"
$ws.Range("F17").Value = "Someone with access to MNs/CNs"

# replica: Method gets explicit call parentheses (still CN-controlled)
$ws.Range("E18").Value = "CNReplication.updateReplicationMetadata()"

# Move the active selection as recorded in the saved view state
[void]$ws.Range("F4").Select()
